# Refresh the "cryptos" price list (Price + Volume(1h) columns, plus a couple
# of coin rows that moved rank and swapped places) to match the latest pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every "Price" cell is stored as text in this sheet (values like "3.249.18" use
# dots as thousands separators, so they are not valid numbers). Force the column
# to Text first so Excel does not silently re-interpret the new values as numbers,
# then drop back to the default "Normal" style once the text is in place so no
# cell ends up with a lingering explicit number format.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "87.932.20"
$ws.Range("E2").Value = "  -0.81%  "

# Row 3
$ws.Range("D3").Value = "3.248.68"
$ws.Range("E3").Value = "  -3.38%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "212.35"
$ws.Range("E5").Value = "  -4.55%  "

# Row 6
$ws.Range("D6").Value = "626.65"
$ws.Range("E6").Value = "  -4.81%  "

# Row 7
$ws.Range("D7").Value = "0.384"
$ws.Range("E7").Value = "  +11.08%  "

# Row 8
$ws.Range("E8").Value = "  +14.77%  "

# Row 9
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.08%  "

# Row 10
$ws.Range("D10").Value = "3.244.96"
$ws.Range("E10").Value = "  -3.45%  "

# Row 11
$ws.Range("D11").Value = "0.574"
$ws.Range("E11").Value = "  -6.54%  "

# Row 12
$ws.Range("D12").Value = "0.188"
$ws.Range("E12").Value = "  +12.39%  "

# Row 13
$ws.Range("E13").Value = "  -3.19%  "

# Row 14
$ws.Range("D14").Value = "5.48"
$ws.Range("E14").Value = "  -0.30%  "

# Row 15
$ws.Range("D15").Value = "34.10"
$ws.Range("E15").Value = "  -3.96%  "

# Row 16
$ws.Range("D16").Value = "3.849.62"
$ws.Range("E16").Value = "  -3.25%  "

# Row 17
$ws.Range("D17").Value = "87.747.53"
$ws.Range("E17").Value = "  -0.25%  "

# Row 18
$ws.Range("D18").Value = "3.250.13"
$ws.Range("E18").Value = "  -3.03%  "

# Row 19
$ws.Range("E19").Value = "  +0.22%  "

# Row 20
$ws.Range("D20").Value = "14.02"
$ws.Range("E20").Value = "  -5.39%  "

# Row 21
$ws.Range("D21").Value = "435.12"
$ws.Range("E21").Value = "  -7.54%  "

# Row 22
$ws.Range("D22").Value = "8.96"
$ws.Range("E22").Value = "  -3.57%  "

# Row 23
$ws.Range("E23").Value = "  -6.98%  "

# Row 24
$ws.Range("E24").Value = "  -1.12%  "

# Row 25
$ws.Range("E25").Value = "  -3.68%  "

# Row 26
$ws.Range("D26").Value = "12.44"
$ws.Range("E26").Value = "  -7.87%  "

# Row 27
$ws.Range("D27").Value = "0.0000142"
$ws.Range("E27").Value = "  +10.14%  "

# Row 28
$ws.Range("D28").Value = "3.420.02"
$ws.Range("E28").Value = "  -2.51%  "

# Row 29
$ws.Range("D29").Value = "77.25"
$ws.Range("E29").Value = "  -3.12%  "

# Row 30
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
$ws.Range("E31").Value = "  -13.61%  "

# Row 32
$ws.Range("E32").Value = "  +0.49%  "

# Row 33
$ws.Range("D33").Value = "569.75"
$ws.Range("E33").Value = "  -5.88%  "

# Row 34
$ws.Range("D34").Value = "8.84"
$ws.Range("E34").Value = "  -6.18%  "

# Row 35
$ws.Range("D35").Value = "1.38"
$ws.Range("E35").Value = "  -11.84%  "

# Row 36
$ws.Range("D36").Value = "7.26"
$ws.Range("E36").Value = "  +3.09%  "

# Row 37
$ws.Range("E37").Value = "  -6.46%  "

# Row 38
$ws.Range("E38").Value = "  -8.93%  "

# Row 39
$ws.Range("D39").Value = "22.87"
$ws.Range("E39").Value = "  -5.64%  "

# Row 40
$ws.Range("E40").Value = "  +5.87%  "

# Row 41
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "21.81"
$ws.Range("E41").Value = "  +0.62%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.08%  "

# Row 43
$ws.Range("D43").Value = "0.402"
$ws.Range("E43").Value = "  -5.45%  "

# Row 44
$ws.Range("E44").Value = "  -6.21%  "

# Row 46
$ws.Range("D46").Value = "151.53"
$ws.Range("E46").Value = "  -4.20%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.135"
$ws.Range("E47").Value = "  +17.86%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "179.67"
$ws.Range("E48").Value = "  -7.21%  "

# Row 49
$ws.Range("E49").Value = "  -5.73%  "

# Row 50
$ws.Range("E50").Value = "  -3.71%  "

# Row 51
$ws.Range("D51").Value = "4.24"
$ws.Range("E51").Value = "  -3.36%  "

# Restore the default cell style now that every Price cell holds literal text.
$ws.Range("D2:D51").Style = "Normal"

